# Refresh crypto price list values (Price & Volume(1h) columns),
# including three coin re-orderings (rows 23/24, 32/33) and one
# coin replacement (row 51: RocketPoolETH -> MultiversX).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price text reads as a plain number (e.g. "1.00", "258.38")
# must be pre-formatted as Text so Excel keeps the literal string (incl.
# trailing zeros) instead of silently coercing it to a numeric value.
$textPriceCells = "D4,D5,D6,D7,D8,D9,D10,D11,D12,D14,D15,D16,D19,D21,D22,D23,D24,D25,D26,D27,D28,D29,D30,D32,D33,D34,D35,D36,D37,D38,D39,D40,D41,D42,D43,D44,D46,D47,D48,D49,D51"
$ws.Range($textPriceCells).NumberFormat = "@"

$ws.Range("D2").Value = "37.422.05"
$ws.Range("E2").Value = "  +2.15%  "
$ws.Range("D3").Value = "2.039.31"
$ws.Range("E3").Value = "  +3.23%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "258.38"
$ws.Range("E5").Value = "  +5.66%  "
$ws.Range("D6").Value = "0.623"
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("D7").Value = "58.38"
$ws.Range("E7").Value = "  -2.79%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "0.391"
$ws.Range("E9").Value = "  +3.07%  "
$ws.Range("D10").Value = "0.0811"
$ws.Range("E10").Value = "  +2.67%  "
$ws.Range("D11").Value = "0.104"
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("D12").Value = "15.08"
$ws.Range("E12").Value = "  +5.87%  "
$ws.Range("D13").Value = "2.326.50"
$ws.Range("E13").Value = "  +2.58%  "
$ws.Range("D14").Value = "0.834"
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("D15").Value = "21.62"
$ws.Range("E15").Value = "  -0.69%  "
$ws.Range("D16").Value = "5.42"
$ws.Range("E16").Value = "  -0.41%  "
$ws.Range("D17").Value = "2.048.48"
$ws.Range("E17").Value = "  +3.36%  "
$ws.Range("D18").Value = "37.380.53"
$ws.Range("E18").Value = "  +2.14%  "
$ws.Range("D19").Value = "70.27"
$ws.Range("E19").Value = "  +0.58%  "
$ws.Range("D20").Value = "0.0₃0864"
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("D21").Value = "5.29"
$ws.Range("E21").Value = "  +3.51%  "
$ws.Range("D22").Value = "229.54"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").Value = "2.65"
$ws.Range("E23").Value = "  +8.18%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").Value = "2.36"
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("D26").Value = "9.29"
$ws.Range("E26").Value = "  +1.09%  "
$ws.Range("D27").Value = "164.00"
$ws.Range("E27").Value = "  +1.35%  "
$ws.Range("D28").Value = "0.139"
$ws.Range("E28").Value = "  -5.23%  "
$ws.Range("D29").Value = "20.03"
$ws.Range("E29").Value = "  +3.35%  "
$ws.Range("D30").Value = "1.37"
$ws.Range("E30").Value = "  +2.11%  "
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "4.82"
$ws.Range("E32").Value = "  +0.38%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "0.0675"
$ws.Range("E33").Value = "  +9.21%  "
$ws.Range("D34").Value = "4.56"
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("D35").Value = "2.50"
$ws.Range("E35").Value = "  +9.60%  "
$ws.Range("D36").Value = "3.57"
$ws.Range("E36").Value = "  +8.15%  "
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("D38").Value = "1.82"
$ws.Range("E38").Value = "  +2.23%  "
$ws.Range("D39").Value = "5.42"
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("D40").Value = "3.02"
$ws.Range("E40").Value = "  +3.75%  "
$ws.Range("D41").Value = "0.0979"
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("D42").Value = "0.0218"
$ws.Range("E42").Value = "  +3.75%  "
$ws.Range("D43").Value = "1.19"
$ws.Range("E43").Value = "  +1.87%  "
$ws.Range("D44").Value = "16.67"
$ws.Range("E44").Value = "  +4.59%  "
$ws.Range("D45").Value = "1.406.25"
$ws.Range("E45").Value = "  +2.99%  "
$ws.Range("D46").Value = "91.89"
$ws.Range("E46").Value = "  +3.05%  "
$ws.Range("D47").Value = "1.06"
$ws.Range("E47").Value = "  +2.91%  "
$ws.Range("D48").Value = "7.50"
$ws.Range("E48").Value = "  +4.16%  "
$ws.Range("D49").Value = "2.09"
$ws.Range("E49").Value = "  +10.97%  "
$ws.Range("E50").Value = "  +2.14%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "45.89"
$ws.Range("E51").Value = "  -0.73%  "
